$d = $word.ActiveDocument

# 1. Append JavaScript-browser sentence to the "Make sure to leave this window open..." paragraph.
$pLeaveWindow = $d.Paragraphs.Item(11)
$pLeaveWindow.Range.InsertAfter(" Make sure your browser has JavaScript enabled, or you may not see a link to the informed consent document below.")

# 2. Add the "*" footnote marker before the trailing colon in the "To access the survey..." paragraph.
$d.Content.Find.Execute(
    "at the bottom of the informed consent document below:", $true, $false, $false, $false, $false,
    $true, 1, $false, "at the bottom of the informed consent document below*:", 2) | Out-Null

# 3. Update the first OSF hyperlink address/display text.
$h1 = $d.Hyperlinks.Item(1)
$h1.Address = "https://osf.io/z95j8/"
$h1.TextToDisplay = "https://osf.io/z95j8/"

# 4. Insert a new paragraph with a second OSF hyperlink right after the first hyperlink paragraph.
$pHyperlink1 = $d.Paragraphs.Item(13)
$pHyperlink1.Range.InsertParagraphAfter()
$pNewHyperlink = $d.Paragraphs.Item(14)
$d.Hyperlinks.Add($pNewHyperlink.Range, "https://osf.io/cjrfm/", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "https://osf.io/cjrfm/") | Out-Null

# 5. Append the randomization explanation paragraph at the very end of the document.
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$pLast.Range.InsertParagraphAfter()
$pNote = $d.Paragraphs.Item($d.Paragraphs.Count)
$pNote.Range.InsertAfter("*Users are randomly presented one of the two survey variants using JavaScript embedded in Mechanical Turk")
